# Bug fixes for Outlier detection
# Expand the adjacency matrix from 7x7 (A1:G7) to 10x10 (A1:J10) and
# update the cell values to the corrected matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1,1,1,1,0,0,0,0,0,0),
    @(1,1,1,0,0,0,0,0,0,0),
    @(1,1,1,0,0,0,0,0,0,0),
    @(1,0,0,1,0,0,0,0,0,0),
    @(0,0,0,0,1,1,1,0,0,0),
    @(0,0,0,0,1,1,0,0,0,0),
    @(0,0,0,0,1,0,1,0,0,0),
    @(0,0,0,0,0,0,0,1,1,0),
    @(0,0,0,0,0,0,0,1,1,1),
    @(0,0,0,0,0,0,0,0,1,1)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
